$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the log entries that are no longer part of this export (rows 5-8)
$ws.Range("A5:G8").EntireRow.Delete()

# Row 2 - existing first log line, refreshed for the 22 Sep deployment
$b2 = $ws.Range("B2")
$b2.Value = "'2023-09-22"
$b2.Style = "Normal"
$ws.Range("C2").Value = "08:36:53"
$ws.Range("D2").Value = "Marc"
$ws.Range("E2").Value = "Geraerts"
$ws.Range("F2").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G2").Value = "Lucp2284"

# Row 3 - updated log entry
$b3 = $ws.Range("B3")
$b3.Value = "'2023-09-22"
$b3.Style = "Normal"
$ws.Range("C3").Value = "09:55:44"
$ws.Range("D3").Value = "test66"
$ws.Range("E3").Value = "test"
$ws.Range("F3").Value = "lucp7894"
$ws.Range("G3").Value = "test66@test.com"

# Row 4 - new (mostly blank) log entry
$b4 = $ws.Range("B4")
$b4.Value = "'2023-09-22"
$b4.Style = "Normal"
$ws.Range("C4").Value = "09:56:06"
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
